$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Through 2022-07-14"
$ws.Range("B1").Value = "July 2022 (through July 14)"

$ws.Range("B2").Value = 8
$ws.Range("P2").Value = 3

$ws.Range("B3").Value = 7
$ws.Range("I3").Value = 3
$ws.Range("AR3").Value = 4

$ws.Range("I8").Value = 5

$ws.Range("I14").Value = 2

$ws.Range("AK15").Value = 2

$ws.Range("B18").Value = 1

$ws.Range("B21").Value = 1

$ws.Range("I36").Value = 2

$ws.Range("B39").Value = 1

$ws.Range("AD41").Value = 1

$ws.Range("W45").Value = 1

$ws.Range("P52").Value = 7

$ws.Range("P53").Value = 5

$ws.Range("AR69").Value = 1

$ws.Range("B97").Value = 1
